$wb = $excel.ActiveWorkbook

# Mapping from the old full Group names (column A) to the new abbreviated
# Group names that now live in the freshly-inserted column C.
function Get-GroupAbbrev($name) {
    if ($name -eq "Roads Places and Environment Group") { return "RPE" }
    if ($name -eq "Rail Group") { return "Rail" }
    if ($name -eq "HSMRPG") { return "HSMRPG" }
    return $name
}

foreach ($sheetName in @("Q1_20_21", "Q4_19_20")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Figure out the last used row on this sheet before we touch anything.
    $lastRow = $ws.Cells.SpecialCells(11).Row

    # Insert a brand-new column before the current column C, which pushes
    # the existing C:K block (NPV .. Benefits Narrative) one column right
    # to D:L.
    $ws.Columns.Item(3).Insert()

    # New header for the inserted column.
    $ws.Range("C2").Value = "Group"

    # Column A held the full "Group" text on data rows (row 3 is the Mars
    # summary row and never had one). Move each value into the new column C
    # as its abbreviation, then clear the old column A cell.
    for ($r = 3; $r -le $lastRow; $r++) {
        $old = $ws.Range("A" + $r).Value2
        if ($old -ne $null -and $old -ne "") {
            $ws.Range("C" + $r).Value = Get-GroupAbbrev($old)
            $ws.Range("A" + $r).Value = $null
        }
    }
}
